$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update the protein file name value in F2 (shared string "Protein_selenium" -> "SYLK003883-PROT.prot")
$ws.Range("F2").Value = "SYLK003883-PROT.prot"

# 2. Give column S (19) an explicit custom width, matching column G's width
$ws.Columns.Item(19).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# 3. Reset row 1 back to the default (auto) row height, removing its explicit 45pt height
$ws.Rows.Item(1).AutoFit()

# 4. Move the active selection from T2 to S1
[void]$ws.Range("S1").Select()
